$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 13, pushing rows 13-23 down to 14-24,
# and carrying per-row formatting (row height) along with them.
$ws.Rows.Item(13).Insert()

# --- Row 10 (Objetivos:) - B/C get the new Portuguese objectives text ---
$ws.Range("B10").Value = 'A reologia é a ciência que estuda o escoamento de materiais. O seu conhecimento é necessário para poder entender o processamento dos materiais poliméricos. A disciplina visa o ensino dos conceitos básicos de reologia de materiais (polímeros fundidos) para o estudante de engenharia de materiais, a disciplina visa também familiarizar o futuro engenheiro com os métodos experimentais para avaliação das propriedades reológicas de materiais poliméricos.Fornecer conhecimentos técnicos para o aluno escolher corretamente a técnica mais adequada de processamento de polímeros, bem como poder manipular e especificar corretamente as matérias primas, os equipamentos de processamento, os moldes e as máquinas adequadamente para determinadas conformações.'
$ws.Range("C10").Value = 'A reologia é a ciência que estuda o escoamento de materiais. O seu conhecimento é necessário para poder entender o processamento dos materiais poliméricos. A disciplina visa o ensino dos conceitos básicos de reologia de materiais (polímeros fundidos) para o estudante de engenharia de materiais, a disciplina visa também familiarizar o futuro engenheiro com os métodos experimentais para avaliação das propriedades reológicas de materiais poliméricos.Fornecer conhecimentos técnicos para o aluno escolher corretamente a técnica mais adequada de processamento de polímeros, bem como poder manipular e especificar corretamente as matérias primas, os equipamentos de processamento, os moldes e as máquinas adequadamente para determinadas conformações.'

# --- Row 13 (new blank row) - gets the "Docentes responsaveis" value ---
$ws.Range("B13").Value = '5840897 - Clodoaldo Saron'
$ws.Range("C13").Value = '5840897 - Clodoaldo Saron'

# --- Row 14 (Programa resumido:) - B/C get the new PT short syllabus ---
$ws.Range("B14").Value = 'Fundamentos de reologia. Processamento de polímeros: matérias-primas, máquinas e moldes.'
$ws.Range("C14").Value = 'Fundamentos de reologia. Processamento de polímeros: matérias-primas, máquinas e moldes.'

# --- Row 16 (Programa:) - B/C get the new PT full syllabus ---
$ws.Range("B16").Value = '1. Introdução a reologia. Tipos de fluxo. 2. Sólidos hookeanos e fluidos newtonianos. 3. Fluidos newtonianos e não newtonianos. 4. Viscoelasticidade. Viscosidade extensional. Diferenças de tensões normais. Variáveis que afetam a viscosidade de polímeros. 5. Importância da Reologia no processamento de polímeros. Fluxos utilizados para caracterizar materiais: fluxo de arraste, fluxos devido a diferença de pressão e escoamento em dutos. 6. Extrusão de polímeros: equipamentos, roscas, matrizes e aplicações. 7. Injeção de polímeros: equipamento, moldes, controle da operação, correção de problemas e aplicações. 8. Outras técnicas de processamento de termoplásticos: sopro, prensagem, termoformagem, calandragem, fiação, rotomoldagem. 9. Blendas e Compósitos Poliméricas: formas de obtenção, miscibilidade, compatibilidade e aplicações. 10. Técnicas de processamento de polímeros termorrígidos: moldagem manual, moldagem por pistola, pultrusão, enrolamento de filamento, prensagem, etc.'
$ws.Range("C16").Value = '1. Introdução a reologia. Tipos de fluxo. 2. Sólidos hookeanos e fluidos newtonianos. 3. Fluidos newtonianos e não newtonianos. 4. Viscoelasticidade. Viscosidade extensional. Diferenças de tensões normais. Variáveis que afetam a viscosidade de polímeros. 5. Importância da Reologia no processamento de polímeros. Fluxos utilizados para caracterizar materiais: fluxo de arraste, fluxos devido a diferença de pressão e escoamento em dutos. 6. Extrusão de polímeros: equipamentos, roscas, matrizes e aplicações. 7. Injeção de polímeros: equipamento, moldes, controle da operação, correção de problemas e aplicações. 8. Outras técnicas de processamento de termoplásticos: sopro, prensagem, termoformagem, calandragem, fiação, rotomoldagem. 9. Blendas e Compósitos Poliméricas: formas de obtenção, miscibilidade, compatibilidade e aplicações. 10. Técnicas de processamento de polímeros termorrígidos: moldagem manual, moldagem por pistola, pultrusão, enrolamento de filamento, prensagem, etc.'

# --- Row 19 (Metodo:) - B/C get the avaliacao text ---
$ws.Range("B19").Value = 'A avaliação será feita por meio de provas escritas.'
$ws.Range("C19").Value = 'A avaliação será feita por meio de provas escritas.'

# --- Row 20 (Criterio:) - B/C get the Nota final formula text ---
$ws.Range("B20").Value = 'A Nota final (NF) será calculada da seguinte maneira: NF = (P1 + P2)/2'
$ws.Range("C20").Value = 'A Nota final (NF) será calculada da seguinte maneira: NF = (P1 + P2)/2'

# --- Row 21 (Norma de recuperacao:) - B/C get the recovery-exam text ---
$ws.Range("B21").Value = 'A recuperação será feita por meio de uma prova escrita (PR) e a média de recuperação (MR) calculada pela fórmula: MR = (NF + PR)/2'
$ws.Range("C21").Value = 'A recuperação será feita por meio de uma prova escrita (PR) e a média de recuperação (MR) calculada pela fórmula: MR = (NF + PR)/2'

# --- Row 22 (Bibliografia:) - B/C get the bibliography text ---
$ws.Range("B22").Value = 'BRETAS, R. E. S.; D´ÁVILA, M. A. Reologia de Polímeros Fundidos, São Carlos, Eduscar, 2005.MANRICH, S. Processamento de termoplásticos – Rosca única, extrusão & matrizes, injeção & moldes,. McCRUM, N. G., BUCKLEY, C. P., BUCKNALl, C. B. Principles of Polymer Engineering, New York, Oxford University Press, 1997.Blass A., Processamento de Polímeros, editora da UFSC.CHAWLA, K. K. Composite Materials Science and Engineering, Spring-Verlag ed., Berlin, 1987.BRETT, A.M.O., BRETT, C.M. Electroquímica: Princípios, métodos e aplicações. Livraria Medina, Coimbra, 1996.FONTANA, M. G. Corrosion Engineering. 3ª Edição. McGraw-Hill, 1987GENTIL, V. Corrosão. 5ª Edição, Rio de Janeiro, Ed. LTC, 2007 RAMANHATAN, L. Corrosão e seu Controle. São Paulo. Ed. Hemus, 1990SHREIR, L.L., JARMAN, R.A., BURSTEIN, G.T. Corrosion. 3ª Edição. Oxford, Butterworth Heinemann, volume 2, 2000WOLYNEC, S. Técnicas Eletroquímicas em Corrosão, EDUSP, São Paulo, 2003'
$ws.Range("C22").Value = 'BRETAS, R. E. S.; D´ÁVILA, M. A. Reologia de Polímeros Fundidos, São Carlos, Eduscar, 2005.MANRICH, S. Processamento de termoplásticos – Rosca única, extrusão & matrizes, injeção & moldes,. McCRUM, N. G., BUCKLEY, C. P., BUCKNALl, C. B. Principles of Polymer Engineering, New York, Oxford University Press, 1997.Blass A., Processamento de Polímeros, editora da UFSC.CHAWLA, K. K. Composite Materials Science and Engineering, Spring-Verlag ed., Berlin, 1987.BRETT, A.M.O., BRETT, C.M. Electroquímica: Princípios, métodos e aplicações. Livraria Medina, Coimbra, 1996.FONTANA, M. G. Corrosion Engineering. 3ª Edição. McGraw-Hill, 1987GENTIL, V. Corrosão. 5ª Edição, Rio de Janeiro, Ed. LTC, 2007 RAMANHATAN, L. Corrosão e seu Controle. São Paulo. Ed. Hemus, 1990SHREIR, L.L., JARMAN, R.A., BURSTEIN, G.T. Corrosion. 3ª Edição. Oxford, Butterworth Heinemann, volume 2, 2000WOLYNEC, S. Técnicas Eletroquímicas em Corrosão, EDUSP, São Paulo, 2003'

